# Update Name of Algo
# Applies updated imputed values (columns A, C, D) produced by the
# RandomForest algorithm run for terrestrial_mammals / combination_1_ABCD
# / ACD / 20 / seed3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C4").Value = -12.3174
$ws.Range("D5").Value = -7.596300000000005
$ws.Range("A8").Value = -22.22470000000002
$ws.Range("D8").Value = -8.523599999999995
$ws.Range("A10").Value = -21.67759999999999
$ws.Range("C11").Value = -13.1931
$ws.Range("A12").Value = -21.67439999999999
$ws.Range("C12").Value = -12.91300000000001
$ws.Range("D12").Value = -8.929800000000004
$ws.Range("D13").Value = -8.854399999999988
$ws.Range("C15").Value = -14.3475
$ws.Range("D15").Value = -8.346299999999998
$ws.Range("C17").Value = -13.3822
$ws.Range("A18").Value = -22.25590000000001
$ws.Range("D21").Value = -7.818799999999991
$ws.Range("A25").Value = -21.7685
$ws.Range("D25").Value = -7.699199999999996
$ws.Range("C26").Value = -12.47510000000001
$ws.Range("C27").Value = -13.37569999999999
$ws.Range("C28").Value = -13.67519999999999
$ws.Range("C32").Value = -13.07270000000001
$ws.Range("D32").Value = -8.913100000000004
$ws.Range("D36").Value = -7.504299999999996
$ws.Range("A37").Value = -19.60649999999999
$ws.Range("C37").Value = -12.97100000000001
$ws.Range("D38").Value = -8.056099999999997
$ws.Range("C41").Value = -12.5754
$ws.Range("D41").Value = -8.316000000000001
$ws.Range("C47").Value = -12.8614
$ws.Range("D50").Value = -7.751899999999998
$ws.Range("C51").Value = -12.084
$ws.Range("D52").Value = -7.883799999999999
$ws.Range("A55").Value = -22.20979999999999
$ws.Range("D59").Value = -8.478700000000002
$ws.Range("C65").Value = -12.10919999999999
$ws.Range("D67").Value = -7.224499999999994
$ws.Range("A68").Value = -21.46349999999999
$ws.Range("C73").Value = -11.78390000000001
$ws.Range("A77").Value = -19.70959999999999
$ws.Range("A78").Value = -19.68809999999998
$ws.Range("A79").Value = -19.74319999999998
$ws.Range("A80").Value = -19.7505
$ws.Range("A81").Value = -22.00790000000001
$ws.Range("A82").Value = -21.91530000000002
$ws.Range("A84").Value = -21.98129999999999
$ws.Range("C84").Value = -13.22789999999999
$ws.Range("D84").Value = -7.971999999999993
$ws.Range("C85").Value = -13.6317
$ws.Range("D86").Value = -8.138299999999997
$ws.Range("D88").Value = -7.718199999999999
$ws.Range("C89").Value = -13.8058
$ws.Range("D89").Value = -8.089899999999998
$ws.Range("C93").Value = -10.3632
$ws.Range("C95").Value = -13.14519999999999
$ws.Range("D95").Value = -7.742399999999996
$ws.Range("C98").Value = -13.02140000000001
$ws.Range("C99").Value = -12.0073
$ws.Range("A101").Value = -20.46159999999998
$ws.Range("C101").Value = -12.26990000000001
$ws.Range("A102").Value = -20.55409999999998
$ws.Range("C102").Value = -12.0853
$ws.Range("D105").Value = -8.157400000000006
